$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are stored as plain text (e.g. "291.34"),
# so force text entry (NumberFormat "@") then restore the default
# "Normal" style so no stray number-format style sticks to the cell.

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "39.944.53"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.216.01"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  -0.04%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "291.34"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "86.66"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("E7").Value = "  -0.55%  "

$ws.Range("E9").Value = "  -1.27%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "30.37"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "50.29"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +5.80%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0778"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("E13").Value = "  +3.01%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "6.42"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "2.561.56"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "13.76"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -2.30%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "2.247.34"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("E18").Value = "  +0.01%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "39.867.75"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +0.14%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "11.06"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -4.36%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.73"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "65.60"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "237.90"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("E27").Value = "  -0.64%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "23.00"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("E30").Value = "  -7.71%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "156.42"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +2.68%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "31.81"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("E33").Value = "  -0.03%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "4.95"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "2.97"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +5.59%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0712"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "2.34"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  +1.31%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "15.21"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -4.72%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "2.093.22"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "3.69"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -2.90%  "

$ws.Range("E44").Value = "  +0.68%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "17.91"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  +1.56%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "9.77"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -2.53%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "1.99"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -8.18%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "2.70"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "2.432.97"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "1.45"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +0.12%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "1.10"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +2.28%  "
